$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Part 1: add <w:noProof/> to the run property of the (first) run carrying
# each inline picture in the five image-only paragraphs.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.InlineShapes.Count -gt 0) {
        $p.Range.Font.NoProofing = 1
    }
}

# ---------------------------------------------------------------------------
# Part 2: append the new "e) MDP ..." section at the end of the document.
# ---------------------------------------------------------------------------

# Step 2a: the existing (empty) final paragraph gets the "e)MDP:Morozov's
# Discrepancy Principle" runs merged into it, preserving its own paragraph
# properties/identity.
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r1 = $d.Range($lastPara.Range.Start, $lastPara.Range.Start)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>e)</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>MDP</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>:Morozov\u2019s</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t xml:space="preserve"> Discrepancy Principle</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$xml1 = $xml1.Replace("\u2019", [string]([char]0x2019))
$r1.InsertXML($xml1)

# Step 2b: append four brand-new paragraphs after it (Disadvantage,
# Advantage, an empty paragraph, and a final paragraph containing a space).
$newLastPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$r2 = $d.Range($newLastPara.Range.End - 1, $newLastPara.Range.End - 1)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>Di</w:t></w:r><w:r><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>sadvantage:MDP</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t xml:space="preserve">(Morozov) needs information about noise variance </w:t></w:r><w:r><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t>whereas CV is purely data-driven.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:lastRenderedPageBreak/><w:t>Advantage: the reconstruction in MDP uses all measurements whereas reconstruction in CV is done using only reconstruction set which is smaller.</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr></w:pPr><w:r><w:rPr><w:sz w:val="48"/><w:szCs w:val="48"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$r2.InsertXML($xml2)
